# Account for different number of bigs and littles
$wb = $excel.ActiveWorkbook

$littles = $wb.Worksheets.Item("Littles")
$pairs   = $wb.Worksheets.Item("Pairs")

# Add a new Little (row 6) to the "Littles" preference sheet.
$littles.Range("A6").Value = "Little 6"
$littles.Range("B6").Value = "Big 3"
$littles.Range("C6").Value = "Big 2"
$littles.Range("D6").Value = "Big 4"
$littles.Range("A7").Select()

# Update the computed pairing sheet: rows 5 & 6 swap, and a new pair is appended.
$pairs.Range("A5").Value = "Big 3"
$pairs.Range("B5").Value = "Little 3"
$pairs.Range("A6").Value = "Big 4"
$pairs.Range("B6").Value = "Little 1"
$pairs.Range("A7").Value = "Big 1*"
$pairs.Range("B7").Value = "Little 6"
